$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.429.03'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.851.19'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.23'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D8").Value = '4.608.23'
$ws.Range("E8").Value = '  +143.80%  '
$ws.Range("D9").Value = '4.729.05'
$ws.Range("E9").Value = '  +118.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07573'
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.2961'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.61'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07726'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.997'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6848'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.99'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.184'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("D19").Value = '29.458.32'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '231.79'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.608'
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.03'
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.405'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("D29").Value = '4.746.18'
$ws.Range("E29").Value = '  +128.99%  '
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05764'
$ws.Range("E31").Value = '  -3.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.257'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.129'
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("E35").Value = '  -1.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.159'
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7167'
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.595'
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("D39").Value = '1.253.76'
$ws.Range("E39").Value = '  +3.53%  '
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("E41").Value = '  +1.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9054'
$ws.Range("E42").Value = '  -0.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.116'
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -0.33%  '
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.138'
$ws.Range("E47").Value = '  -3.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.204'
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4022'
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.684'
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1124'
$ws.Range("E51").Value = '  -0.60%  '
